$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = '100%>gen/cap>95%'
$ws.Range('B2').Value = 0
$ws.Range('C2:G2').Clear()
$ws.Range('H2').Value = 0.07699423965621355

$ws.Range('A3').Value = '105%>gen/cap>100%'
$ws.Range('B3').Value = 0.002136065079397322
$ws.Range('C3').Value = 0.000714064924072249
$ws.Range('D3').Value = 1.442321740502984
$ws.Range('E3').Value = 0.001357445521550493
$ws.Range('F3').Value = 0.0007365197392141638
$ws.Range('G3').Value = 0.003535610419580481
$ws.Range('H3').Value = 0.07913030473561088

$ws.Range('A4').Value = '110%>gen/cap>105%'
$ws.Range('B4').Value = 0.01159257279452226
$ws.Range('C4').Value = 0.001219869276552708
$ws.Range('D4').Value = 8.283266409881842
$ws.Range('E4').Value = 0.000000002455491085022402
$ws.Range('F4').Value = 0.009201666129321493
$ws.Range('G4').Value = 0.01398347945972302
$ws.Range('H4').Value = 0.08858681245073581

$ws.Range('A5').Value = '115%>gen/cap>110%'
$ws.Range('B5').Value = 0.0531000885534139
$ws.Range('C5').Value = 0.002471498364642754
$ws.Range('D5').Value = 21.81244191468928
$ws.Range('E5').Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000003614122183643611
$ws.Range('F5').Value = 0.04825602623055938
$ws.Range('G5').Value = 0.05794415087626842
$ws.Range('H5').Value = 0.1300943282096275

$ws.Range('A6').Value = '120%>gen/cap>115%'
$ws.Range('B6').Value = 0.06627064596797026
$ws.Range('C6').Value = 0.003612342372319574
$ws.Range('D6').Value = 19.51825469578559
$ws.Range('E6').Value = 0.00000000000000000000000001907706688816555
$ws.Range('F6').Value = 0.05919056869522173
$ws.Range('G6').Value = 0.0733507232407188
$ws.Range('H6').Value = 0.1432648856241838

$ws.Range('A7').Value = '125%>gen/cap>120%'
$ws.Range('B7').Value = 0.05303816425748428
$ws.Range('C7').Value = 0.0154542645840133
$ws.Range('D7').Value = 18.22876526854494
$ws.Range('E7').Value = 0.2749095327643388
$ws.Range('F7').Value = 0.02274831094093235
$ws.Range('G7').Value = 0.08332801757403621
$ws.Range('H7').Value = 0.1300324039136978

$ws.Range('A8').Value = '130%>gen/cap>125%'
$ws.Range('B8').Value = 0.06357106548722641
$ws.Range('C8').Value = 0.01182015119928298
$ws.Range('D8').Value = 14.9415281229764
$ws.Range('E8').Value = 0.09720988479359305
$ws.Range('F8').Value = 0.04040395282381733
$ws.Range('G8').Value = 0.0867381781506355
$ws.Range('H8').Value = 0.14056530514344

$ws.Range('A9').Value = '135%>gen/cap>130%'
$ws.Range('B9').Value = 0.05861274728887465
$ws.Range('C9').Value = 0.0050077386367045
$ws.Range('D9').Value = 17.80121790097613
$ws.Range('E9').Value = 0.02970079211986037
$ws.Range('F9').Value = 0.04879773942199888
$ws.Range('G9').Value = 0.06842775515575042
$ws.Range('H9').Value = 0.1356069869450882

$ws.Range('A10').Value = '20%>gen/cap'
$ws.Range('B10').Value = -0.07699423965621355
$ws.Range('C10').Value = 0.0004765546690353926
$ws.Range('D10').Value = -174.5266124724704
$ws.Range('E10').Value = 0
$ws.Range('F10').Value = -0.07792827213294053
$ws.Range('G10').Value = -0.07606020717948657
$ws.Range('H10').Value = 0

$ws.Range('A11').Value = '25%>gen/cap>20%'
$ws.Range('B11').Value = -0.02888212564271414
$ws.Range('C11').Value = 0.0004980814961774719
$ws.Range('D11').Value = -61.68005171812727
$ws.Range('E11').Value = 0
$ws.Range('F11').Value = -0.02985835004141893
$ws.Range('G11').Value = -0.02790590124400934
$ws.Range('H11').Value = 0.04811211401349941

$ws.Range('A12').Value = '30%>gen/cap>25%'
$ws.Range('B12').Value = -0.02402224049891178
$ws.Range('C12').Value = 0.0004883737916974052
$ws.Range('D12').Value = -52.23973010935477
$ws.Range('E12').Value = 0
$ws.Range('F12').Value = -0.02497943809540563
$ws.Range('G12').Value = -0.02306504290241793
$ws.Range('H12').Value = 0.05297199915730177

$ws.Range('A13').Value = '35%>gen/cap>30%'
$ws.Range('B13').Value = -0.02020633602369902
$ws.Range('C13').Value = 0.0004846976399385352
$ws.Range('D13').Value = -42.88596983142469
$ws.Range('E13').Value = 0
$ws.Range('F13').Value = -0.02115632847994901
$ws.Range('G13').Value = -0.01925634356744903
$ws.Range('H13').Value = 0.05678790363251453

$ws.Range('A14').Value = '40%>gen/cap>35%'
$ws.Range('B14').Value = -0.01749104368499764
$ws.Range('C14').Value = 0.0004731942471222773
$ws.Range('D14').Value = -37.93460796680762
$ws.Range('E14').Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000009361632976665488
$ws.Range('F14').Value = -0.01841848984119145
$ws.Range('G14').Value = -0.01656359752880383
$ws.Range('H14').Value = 0.05950319597121591

$ws.Range('A15').Value = '45%>gen/cap>40%'
$ws.Range('B15').Value = -0.01477166331715771
$ws.Range('C15').Value = 0.0004719296174076465
$ws.Range('D15').Value = -33.50377036548491
$ws.Range('E15').Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000005351999412689911
$ws.Range('F15').Value = -0.01569663085005864
$ws.Range('G15').Value = -0.01384669578425677
$ws.Range('H15').Value = 0.06222257633905584

$ws.Range('A16').Value = '50%>gen/cap>45%'
$ws.Range('B16').Value = -0.0109059500294208
$ws.Range('C16').Value = 0.0004735548231036322
$ws.Range('D16').Value = -24.60616739008181
$ws.Range('E16').Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000007461792989506166
$ws.Range('F16').Value = -0.0118341029308993
$ws.Range('G16').Value = -0.009977797127942304
$ws.Range('H16').Value = 0.06608828962679275

$ws.Range('A17').Value = '55%>gen/cap>50%'
$ws.Range('B17').Value = -0.01128488307986634
$ws.Range('C17').Value = 0.0004882389048789513
$ws.Range('D17').Value = -23.18263302348532
$ws.Range('E17').Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000007117298691780633
$ws.Range('F17').Value = -0.01224181632010759
$ws.Range('G17').Value = -0.01032794983962509
$ws.Range('H17').Value = 0.0657093565763472

$ws.Range('A18').Value = '60%>gen/cap>55%'
$ws.Range('B18').Value = -0.007760298307205492
$ws.Range('C18').Value = 0.0004969206605491901
$ws.Range('D18').Value = -14.61443822485369
$ws.Range('E18').Value = 0.00000000000000000000000000001618521401045535
$ws.Range('F18').Value = -0.008734247508980778
$ws.Range('G18').Value = -0.006786349105430205
$ws.Range('H18').Value = 0.06923394134900805

$ws.Range('A19').Value = '65%>gen/cap>60%'
$ws.Range('B19').Value = -0.007789544702405066
$ws.Range('C19').Value = 0.0004937492989384938
$ws.Range('D19').Value = -16.18269818329463
$ws.Range('E19').Value = 0.0000000000000000000000000000000000000000000000000188988338419006
$ws.Range('F19').Value = -0.00875727814952183
$ws.Range('G19').Value = -0.006821811255288304
$ws.Range('H19').Value = 0.06920469495380849

$ws.Range('A20').Value = '70%>gen/cap>65%'
$ws.Range('B20').Value = -0.007957848445520917
$ws.Range('C20').Value = 0.0005071988306013597
$ws.Range('D20').Value = -16.34893509397246
$ws.Range('E20').Value = 0.0000000000000000000000000000000000000000000001646061341705842
$ws.Range('F20').Value = -0.00895194256335215
$ws.Range('G20').Value = -0.006963754327689683
$ws.Range('H20').Value = 0.06903639121069263

$ws.Range('A21').Value = '75%>gen/cap>70%'
$ws.Range('B21').Value = -0.006266542781214153
$ws.Range('C21').Value = 0.0005189023264325247
$ws.Range('D21').Value = -11.83704454276323
$ws.Range('E21').Value = 0.00000000000000000000000000008544350706639894
$ws.Range('F21').Value = -0.007283575387421551
$ws.Range('G21').Value = -0.005249510175006755
$ws.Range('H21').Value = 0.0707276968749994

$ws.Range('A22').Value = '80%>gen/cap>75%'
$ws.Range('B22').Value = -0.004002812593140008
$ws.Range('C22').Value = 0.0005069924425839406
$ws.Range('D22').Value = -7.877349868079113
$ws.Range('E22').Value = 0.000000000000003783802561011822
$ws.Range('F22').Value = -0.004996502165609841
$ws.Range('G22').Value = -0.003009123020670176
$ws.Range('H22').Value = 0.07299142706307354

$ws.Range('A23').Value = '85%>gen/cap>80%'
$ws.Range('B23').Value = -0.00242980919556134
$ws.Range('C23').Value = 0.000503812841106655
$ws.Range('D23').Value = -4.493664753797132
$ws.Range('E23').Value = 0.0006715006414980336
$ws.Range('F23').Value = -0.003417266825376111
$ws.Range('G23').Value = -0.00144235156574657
$ws.Range('H23').Value = 0.07456443046065221

$ws.Range('A24').Value = '90%>gen/cap>85%'
$ws.Range('B24').Value = -0.001303467205928907
$ws.Range('C24').Value = 0.0005209313027690822
$ws.Range('D24').Value = -3.122159773580741
$ws.Range('E24').Value = 0.1354620591317127
$ws.Range('F24').Value = -0.002324476537535025
$ws.Range('G24').Value = -0.0002824578743227876
$ws.Range('H24').Value = 0.07569077245028465

$ws.Range('A25').Value = '95%>gen/cap>90%'
$ws.Range('B25').Value = -0.002581400184398563
$ws.Range('C25').Value = 0.0005107539393873822
$ws.Range('D25').Value = -5.817468755401923
$ws.Range('E25').Value = 0.0001595666129936088
$ws.Range('F25').Value = -0.003582462229110675
$ws.Range('G25').Value = -0.001580338139686451
$ws.Range('H25').Value = 0.07441283947181498

$ws.Range('A26').Value = 'gen/cap>135%'
$ws.Range('B26').Value = 0.06388423309373906
$ws.Range('C26').Value = 0.001596768634045467
$ws.Range('D26').Value = 35.52484847518591
$ws.Range('E26').Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000006967042354835999
$ws.Range('F26').Value = 0.06075461536233442
$ws.Range('G26').Value = 0.06701385082514369
$ws.Range('H26').Value = 0.1408784727499526
